$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Metadata")

# Update Date row (row 9, column B)
$ws.Range("B9").Value = "2025-07-11T12:29:53+00:00"

# Update Jurisdiction row (row 12, column B)
$ws.Range("B12").Value = "FRANCE"
